# Update "想去人数" (want-to-go count) figures for a few 漫展 (convention) entries.
# These values live on both the "展览" sheet and the consolidated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 544
$ws1.Range("F8").Value = 2311
$ws1.Range("F10").Value = 5744

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 544
$ws4.Range("F11").Value = 2311
$ws4.Range("F13").Value = 5744
